$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "CompFwdRC"
$ws.Range("B8").Value = 91
$ws.Range("C8").Value = 273

$ws.Range("F8").Select()
